# chore: update Sheets via scheduled runner
# Refresh the cached marketboard-driven profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) on each job's Leve-profit sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 15866.333
$ws.Range("I64").Value = 12599
$ws.Range("J64").Value = 17500
$ws.Range("K64").Value = 12599
$ws.Range("L64").Value = 17500
$ws.Range("M64").Value = -12351
$ws.Range("N64").Value = -17996
$ws.Range("H67").Value = 15866.333
$ws.Range("I67").Value = 12599
$ws.Range("J67").Value = 17500
$ws.Range("K67").Value = 12599
$ws.Range("L67").Value = 17500
$ws.Range("M67").Value = -11741
$ws.Range("N67").Value = -19216
$ws.Range("H88").Value = 7096
$ws.Range("J88").Value = 7096
$ws.Range("L88").Value = 7096
$ws.Range("N88").Value = -7908
$ws.Range("H91").Value = 7096
$ws.Range("J91").Value = 7096
$ws.Range("L91").Value = 7096
$ws.Range("N91").Value = -9904
$ws.Range("H92").Value = 5040.96
$ws.Range("I92").Value = 5425
$ws.Range("J92").Value = 624.5
$ws.Range("K92").Value = 5425
$ws.Range("L92").Value = 624.5
$ws.Range("M92").Value = -4177
$ws.Range("N92").Value = -3120.5
$ws.Range("H107").Value = 759.05884
$ws.Range("I107").Value = 759.05884
$ws.Range("K107").Value = 759.05884
$ws.Range("M107").Value = 1160.94116
$ws.Range("H113").Value = 72642.94
$ws.Range("I113").Value = 155711.72
$ws.Range("J113").Value = 14494.8
$ws.Range("K113").Value = 155711.72
$ws.Range("L113").Value = 14494.8
$ws.Range("M113").Value = -152457.72
$ws.Range("N113").Value = -21002.8
$ws.Range("H135").Value = 1402.9286
$ws.Range("I135").Value = 583.375
$ws.Range("J135").Value = 2495.6667
$ws.Range("K135").Value = 5250.375
$ws.Range("L135").Value = 22461.0003
$ws.Range("M135").Value = -2715.375
$ws.Range("N135").Value = -27531.0003
$ws.Range("H136").Value = 80780
$ws.Range("J136").Value = 80780
$ws.Range("L136").Value = 80780
$ws.Range("N136").Value = -90980
$ws.Range("H141").Value = 1521.9286
$ws.Range("I141").Value = 1464.6
$ws.Range("J141").Value = 1999.6666
$ws.Range("K141").Value = 4393.799999999999
$ws.Range("L141").Value = 5998.9998
$ws.Range("M141").Value = 786.2000000000007
$ws.Range("N141").Value = -16358.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5361.4
$ws.Range("I61").Value = 5097.4165
$ws.Range("K61").Value = 5097.4165
$ws.Range("M61").Value = -4885.4165
$ws.Range("H74").Value = 9551.727999999999
$ws.Range("I74").Value = 1345.9166
$ws.Range("K74").Value = 1345.9166
$ws.Range("M74").Value = -471.9166
$ws.Range("H77").Value = 9551.727999999999
$ws.Range("I77").Value = 1345.9166
$ws.Range("K77").Value = 6729.583000000001
$ws.Range("M77").Value = -2361.583000000001
$ws.Range("H110").Value = 8290.941000000001
$ws.Range("I110").Value = 12237.125
$ws.Range("K110").Value = 12237.125
$ws.Range("M110").Value = -10192.125
$ws.Range("H122").Value = 1410.75
$ws.Range("I122").Value = 1274.4736
$ws.Range("K122").Value = 3823.4208
$ws.Range("M122").Value = -1373.4208
$ws.Range("H136").Value = 5361.4
$ws.Range("I136").Value = 5097.4165
$ws.Range("K136").Value = 15292.2495
$ws.Range("M136").Value = -12742.2495

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 10050
$ws.Range("I95").Value = 4567
$ws.Range("J95").Value = 10548.454
$ws.Range("K95").Value = 4567
$ws.Range("L95").Value = 10548.454
$ws.Range("M95").Value = -1821
$ws.Range("N95").Value = -16040.454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 218313.6
$ws.Range("J9").Value = 218313.6
$ws.Range("L9").Value = 218313.6
$ws.Range("N9").Value = -218649.6
$ws.Range("H28").Value = 8326.546
$ws.Range("J28").Value = 8059.3
$ws.Range("L28").Value = 8059.3
$ws.Range("N28").Value = -8549.299999999999
$ws.Range("H31").Value = 42431.54
$ws.Range("I31").Value = 64129.438
$ws.Range("K31").Value = 64129.438
$ws.Range("M31").Value = -63834.438
$ws.Range("H34").Value = 42431.54
$ws.Range("I34").Value = 64129.438
$ws.Range("K34").Value = 64129.438
$ws.Range("M34").Value = -63927.438
$ws.Range("H58").Value = 3626.5
$ws.Range("I58").Value = 3287.4285
$ws.Range("K58").Value = 3287.4285
$ws.Range("M58").Value = -3084.4285
$ws.Range("H122").Value = 1863.3334
$ws.Range("I122").Value = 1895
$ws.Range("K122").Value = 5685
$ws.Range("M122").Value = -3235
$ws.Range("H134").Value = 27933
$ws.Range("I134").Value = 16215.134
$ws.Range("K134").Value = 48645.402
$ws.Range("M134").Value = -46110.402
$ws.Range("H136").Value = 3626.5
$ws.Range("I136").Value = 3287.4285
$ws.Range("K136").Value = 9862.2855
$ws.Range("M136").Value = -7312.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2563.5
$ws.Range("J122").Value = 3284.6667
$ws.Range("L122").Value = 29562.0003
$ws.Range("N122").Value = -34462.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 20000000
$ws.Range("I11").Value = 20000000
$ws.Range("K11").Value = 20000000
$ws.Range("M11").Value = -19999861
$ws.Range("H95").Value = 15499.667
$ws.Range("J95").Value = 15499.667
$ws.Range("L95").Value = 15499.667
$ws.Range("N95").Value = -20991.667
$ws.Range("H102").Value = 2857.6924
$ws.Range("I102").Value = 2845.8333
$ws.Range("K102").Value = 2845.8333
$ws.Range("M102").Value = -1223.8333
$ws.Range("H113").Value = 1703.6364
$ws.Range("I113").Value = 1875
$ws.Range("J113").Value = 1246.6666
$ws.Range("K113").Value = 1875
$ws.Range("L113").Value = 1246.6666
$ws.Range("M113").Value = 295
$ws.Range("N113").Value = -5586.6666
$ws.Range("H122").Value = 1682.8334
$ws.Range("I122").Value = 1049.25
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 3147.75
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = -697.75
$ws.Range("N122").Value = -13750
$ws.Range("H132").Value = 4437.1113
$ws.Range("J132").Value = 6750
$ws.Range("L132").Value = 20250
$ws.Range("N132").Value = -25310
$ws.Range("H138").Value = 59500
$ws.Range("J138").Value = 59500
$ws.Range("L138").Value = 59500
$ws.Range("N138").Value = -69780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 252.96153
$ws.Range("J55").Value = 354.14285
$ws.Range("L55").Value = 354.14285
$ws.Range("N55").Value = -700.14285
$ws.Range("H132").Value = 2795.2903
$ws.Range("I132").Value = 2367.3794
$ws.Range("K132").Value = 7102.138199999999
$ws.Range("M132").Value = -4572.138199999999
$ws.Range("H136").Value = 5172.0435
$ws.Range("I136").Value = 4689.8887
$ws.Range("J136").Value = 6907.8
$ws.Range("K136").Value = 14069.6661
$ws.Range("L136").Value = 20723.4
$ws.Range("M136").Value = -11519.6661
$ws.Range("N136").Value = -25823.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1666.4445
$ws.Range("J96").Value = 1749.75
$ws.Range("L96").Value = 1749.75
$ws.Range("N96").Value = -4495.75
$ws.Range("H97").Value = 35286
$ws.Range("J97").Value = 35286
$ws.Range("L97").Value = 35286
$ws.Range("N97").Value = -37268
$ws.Range("H132").Value = 4075.1562
$ws.Range("I132").Value = 3748.3845
$ws.Range("K132").Value = 11245.1535
$ws.Range("M132").Value = -8715.1535
$ws.Range("H136").Value = 2920.6875
$ws.Range("I136").Value = 3544
$ws.Range("K136").Value = 10632
$ws.Range("M136").Value = -8082
